$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1749153331348/OwQjboHWXD-camera-video-segment-driver-1749153333848.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSB4X3ZMX2%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T190032Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDxgYTC1nnpBZJ7WqrtkCyMeteVbbymv5IXJvuF%2FzT7JgIgdc8nVtRd%2BsKiNJYMt3EDs6HplPzi8eKnJliHFH4qUAIq3QMIZBAEGgw3ODEyMDQ5NDIyNDQiDAzID9Rm3z6ueshvlCq6A23WAQUt1e02bJKfvwxg7I0N33WRc%2BhHyjUNW5KH%2BQY6mbsk38VyCBB%2B%2BussIItjluiMNHdEB82lrBDZ8peURaHN8VY6sPO4v5eK0pT6beuaQ%2FI0hJZzsmownYoEqV5NBw5JPEEfxEc3LvZb0TAszGBnrFh8cKaTjZ80eZZ4JSKTB54cDw4iiFEO64JIP9MEhl5hMwsy69Y3eo60wk4e0VGy9T2y%2BliodjafP3vFKHGvugoqHnvSGS%2BYxXI%2FIEe131IvAidLHkW7oUaVhJKLUnESFkUW9A0Baj0En8dKz%2BrvnAw6gVelEbMG0eGG91UakF%2BENuPV%2Bdla71%2BFDYRhmAVL6tNPbv66%2B9i%2FGuX%2F%2BOsdND8wdTV6cAutqXa4DiYqIQIdpe7qyJC9qD8xrwXWITjM94S0aRcxAf4Kyd%2FCQBdQx5fwNsvrGDwIFnGcewuez514HpbFl0vGPdXZR3H2kLgnq5e7q%2F%2FHBeFbhco%2BIYuLySmvCj%2BQvDhDsDvzYPQtYEnlAuipBfXZeD9C3l3YFy7IE599msXvT6bkmQrYXk%2BDpGV4t1Xdwv7n5cb4OgRy0rclUZFSPTbFG0Awv%2B2MwgY6pQGAK7vQWLHCuG6P%2BS9Vnlmx1N28g8ieynyLkLP%2FMZei7IrAz23UjgdkVVv1x1roM167oEF106Mr0W6sbyoFifEjlm195oZ0G1k52VQNuvHAaRueakVId%2ByETAJVHXBc9Tja8SzgpzYo1P%2FdhTaqVebmWqJbNcWxZKQMCbPXZPPUesflsG8LeKP9ymizrqPdWtFUx%2FdNZvAOCAriNkEslV0E7wM6h5I%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2003%3A00%3A32%20GMT&X-Amz-Signature=cc7c5bb4d1e17d716ec47e5d81938b125ed30358d35abd09b13339bc34a85fb9'

$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1749153331348/oV5nVwJUOa-camera-video-segment-1749153333848.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSB4X3ZMX2%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T190032Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDxgYTC1nnpBZJ7WqrtkCyMeteVbbymv5IXJvuF%2FzT7JgIgdc8nVtRd%2BsKiNJYMt3EDs6HplPzi8eKnJliHFH4qUAIq3QMIZBAEGgw3ODEyMDQ5NDIyNDQiDAzID9Rm3z6ueshvlCq6A23WAQUt1e02bJKfvwxg7I0N33WRc%2BhHyjUNW5KH%2BQY6mbsk38VyCBB%2B%2BussIItjluiMNHdEB82lrBDZ8peURaHN8VY6sPO4v5eK0pT6beuaQ%2FI0hJZzsmownYoEqV5NBw5JPEEfxEc3LvZb0TAszGBnrFh8cKaTjZ80eZZ4JSKTB54cDw4iiFEO64JIP9MEhl5hMwsy69Y3eo60wk4e0VGy9T2y%2BliodjafP3vFKHGvugoqHnvSGS%2BYxXI%2FIEe131IvAidLHkW7oUaVhJKLUnESFkUW9A0Baj0En8dKz%2BrvnAw6gVelEbMG0eGG91UakF%2BENuPV%2Bdla71%2BFDYRhmAVL6tNPbv66%2B9i%2FGuX%2F%2BOsdND8wdTV6cAutqXa4DiYqIQIdpe7qyJC9qD8xrwXWITjM94S0aRcxAf4Kyd%2FCQBdQx5fwNsvrGDwIFnGcewuez514HpbFl0vGPdXZR3H2kLgnq5e7q%2F%2FHBeFbhco%2BIYuLySmvCj%2BQvDhDsDvzYPQtYEnlAuipBfXZeD9C3l3YFy7IE599msXvT6bkmQrYXk%2BDpGV4t1Xdwv7n5cb4OgRy0rclUZFSPTbFG0Awv%2B2MwgY6pQGAK7vQWLHCuG6P%2BS9Vnlmx1N28g8ieynyLkLP%2FMZei7IrAz23UjgdkVVv1x1roM167oEF106Mr0W6sbyoFifEjlm195oZ0G1k52VQNuvHAaRueakVId%2ByETAJVHXBc9Tja8SzgpzYo1P%2FdhTaqVebmWqJbNcWxZKQMCbPXZPPUesflsG8LeKP9ymizrqPdWtFUx%2FdNZvAOCAriNkEslV0E7wM6h5I%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2003%3A00%3A32%20GMT&X-Amz-Signature=3220a093aa35a6bf719419f3e9fda54becea389eafbf99c2577896cfe90ca0ed'

$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109254/1749147387468/GKgOu1ZeR3-camera-video-segment-driver-1749147389968.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSB4X3ZMX2%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T190031Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDxgYTC1nnpBZJ7WqrtkCyMeteVbbymv5IXJvuF%2FzT7JgIgdc8nVtRd%2BsKiNJYMt3EDs6HplPzi8eKnJliHFH4qUAIq3QMIZBAEGgw3ODEyMDQ5NDIyNDQiDAzID9Rm3z6ueshvlCq6A23WAQUt1e02bJKfvwxg7I0N33WRc%2BhHyjUNW5KH%2BQY6mbsk38VyCBB%2B%2BussIItjluiMNHdEB82lrBDZ8peURaHN8VY6sPO4v5eK0pT6beuaQ%2FI0hJZzsmownYoEqV5NBw5JPEEfxEc3LvZb0TAszGBnrFh8cKaTjZ80eZZ4JSKTB54cDw4iiFEO64JIP9MEhl5hMwsy69Y3eo60wk4e0VGy9T2y%2BliodjafP3vFKHGvugoqHnvSGS%2BYxXI%2FIEe131IvAidLHkW7oUaVhJKLUnESFkUW9A0Baj0En8dKz%2BrvnAw6gVelEbMG0eGG91UakF%2BENuPV%2Bdla71%2BFDYRhmAVL6tNPbv66%2B9i%2FGuX%2F%2BOsdND8wdTV6cAutqXa4DiYqIQIdpe7qyJC9qD8xrwXWITjM94S0aRcxAf4Kyd%2FCQBdQx5fwNsvrGDwIFnGcewuez514HpbFl0vGPdXZR3H2kLgnq5e7q%2F%2FHBeFbhco%2BIYuLySmvCj%2BQvDhDsDvzYPQtYEnlAuipBfXZeD9C3l3YFy7IE599msXvT6bkmQrYXk%2BDpGV4t1Xdwv7n5cb4OgRy0rclUZFSPTbFG0Awv%2B2MwgY6pQGAK7vQWLHCuG6P%2BS9Vnlmx1N28g8ieynyLkLP%2FMZei7IrAz23UjgdkVVv1x1roM167oEF106Mr0W6sbyoFifEjlm195oZ0G1k52VQNuvHAaRueakVId%2ByETAJVHXBc9Tja8SzgpzYo1P%2FdhTaqVebmWqJbNcWxZKQMCbPXZPPUesflsG8LeKP9ymizrqPdWtFUx%2FdNZvAOCAriNkEslV0E7wM6h5I%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2003%3A00%3A31%20GMT&X-Amz-Signature=a02506321b0429eeebe67b4e6e3cd3b8a43e4314a14e6056ca1350ea2345f6f5'

$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109254/1749147387468/dQoeD1EUOr-camera-video-segment-1749147389968.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSB4X3ZMX2%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T190031Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDxgYTC1nnpBZJ7WqrtkCyMeteVbbymv5IXJvuF%2FzT7JgIgdc8nVtRd%2BsKiNJYMt3EDs6HplPzi8eKnJliHFH4qUAIq3QMIZBAEGgw3ODEyMDQ5NDIyNDQiDAzID9Rm3z6ueshvlCq6A23WAQUt1e02bJKfvwxg7I0N33WRc%2BhHyjUNW5KH%2BQY6mbsk38VyCBB%2B%2BussIItjluiMNHdEB82lrBDZ8peURaHN8VY6sPO4v5eK0pT6beuaQ%2FI0hJZzsmownYoEqV5NBw5JPEEfxEc3LvZb0TAszGBnrFh8cKaTjZ80eZZ4JSKTB54cDw4iiFEO64JIP9MEhl5hMwsy69Y3eo60wk4e0VGy9T2y%2BliodjafP3vFKHGvugoqHnvSGS%2BYxXI%2FIEe131IvAidLHkW7oUaVhJKLUnESFkUW9A0Baj0En8dKz%2BrvnAw6gVelEbMG0eGG91UakF%2BENuPV%2Bdla71%2BFDYRhmAVL6tNPbv66%2B9i%2FGuX%2F%2BOsdND8wdTV6cAutqXa4DiYqIQIdpe7qyJC9qD8xrwXWITjM94S0aRcxAf4Kyd%2FCQBdQx5fwNsvrGDwIFnGcewuez514HpbFl0vGPdXZR3H2kLgnq5e7q%2F%2FHBeFbhco%2BIYuLySmvCj%2BQvDhDsDvzYPQtYEnlAuipBfXZeD9C3l3YFy7IE599msXvT6bkmQrYXk%2BDpGV4t1Xdwv7n5cb4OgRy0rclUZFSPTbFG0Awv%2B2MwgY6pQGAK7vQWLHCuG6P%2BS9Vnlmx1N28g8ieynyLkLP%2FMZei7IrAz23UjgdkVVv1x1roM167oEF106Mr0W6sbyoFifEjlm195oZ0G1k52VQNuvHAaRueakVId%2ByETAJVHXBc9Tja8SzgpzYo1P%2FdhTaqVebmWqJbNcWxZKQMCbPXZPPUesflsG8LeKP9ymizrqPdWtFUx%2FdNZvAOCAriNkEslV0E7wM6h5I%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2003%3A00%3A31%20GMT&X-Amz-Signature=1e2e6bf8977505b4b684bd68f8fbb055a00c3a9bdab9eac10853dd7eb8d9d968'
